$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.419.31"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "2.070.57"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "235.14"
$ws.Range("E5").Value = "  -0.86%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "57.45"
$ws.Range("E8").Value = "  -1.29%  "

$ws.Range("E9").Value = "  +3.31%  "

$ws.Range("D10").Value = "0.0773"
$ws.Range("E10").Value = "  +1.46%  "

$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").Value = "2.377.05"
$ws.Range("E12").Value = "  +0.41%  "

$ws.Range("D13").Value = "14.45"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("D14").Value = "20.74"
$ws.Range("E14").Value = "  -1.57%  "

$ws.Range("D15").Value = "0.779"
$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").Value = "2.072.48"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").Value = "37.372.01"
$ws.Range("E18").Value = "  -0.54%  "

$ws.Range("D19").Value = "6.37"
$ws.Range("E19").Value = "  +3.74%  "

$ws.Range("D20").Value = "69.73"
$ws.Range("E20").Value = "  +0.75%  "

$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  +0.38%  "

$ws.Range("D22").Value = "227.02"
$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("E25").Value = "  -1.81%  "

$ws.Range("D26").Value = "167.05"
$ws.Range("E26").Value = "  +1.89%  "

$ws.Range("D27").Value = "8.85"
$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").Value = "1.43"
$ws.Range("E28").Value = "  -3.07%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "19.15"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.127"
$ws.Range("E30").Value = "  +0.81%  "

$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("D32").Value = "4.54"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").Value = "0.0619"
$ws.Range("E33").Value = "  -1.14%  "

$ws.Range("D34").Value = "4.56"
$ws.Range("E34").Value = "  +1.59%  "

$ws.Range("E35").Value = "  -2.88%  "

$ws.Range("D36").Value = "1.78"
$ws.Range("E36").Value = "  -0.25%  "

$ws.Range("D37").Value = "3.31"
$ws.Range("E37").Value = "  -1.80%  "

$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").Value = "5.67"
$ws.Range("E39").Value = "  -3.48%  "

$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D40").Value = "0.0969"
$ws.Range("E40").Value = "  -2.43%  "

$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "2.95"
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("D42").Value = "98.37"
$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("D43").Value = "1.482.48"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("D45").Value = "1.17"
$ws.Range("E45").Value = "  +0.93%  "

$ws.Range("D46").Value = "4.11"
$ws.Range("E46").Value = "  -8.70%  "

$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("D48").Value = "15.30"
$ws.Range("E48").Value = "  -3.89%  "

$ws.Range("D49").Value = "7.25"
$ws.Range("E49").Value = "  +0.57%  "

$ws.Range("D50").Value = "2.96"
$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("D51").Value = "2.264.36"
$ws.Range("E51").Value = "  +0.40%  "
